$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")

# Version: 2.0.0-sd-202312-matchbox-patch -> 2.0.0-sd-202406-matchbox-patch
$meta.Range("B3").Value = "2.0.0-sd-202406-matchbox-patch"

# Date: 2024-03-12T18:28:21+01:00 -> 2024-06-19T17:47:42+02:00
$meta.Range("B8").Value = "2024-06-19T17:47:42+02:00"

# Contact: "No display for ContactDetail" -> HL7 Structured Documents contact
$meta.Range("B10").Value = "HL7 International - Structured Documents (http://www.hl7.org/Special/committees/structure, structdog@lists.HL7.org)"

# --- Elements sheet updates ---
$elements = $wb.Worksheets.Item("Elements")

# Binding Value Set for the SXCM_TS.operator row (row 5): update ValueSet URL
$elements.Range("Z5").Value = "http://hl7.org/cda/stds/core/ValueSet/CDASetOperator"

# Column Z width grew (best-fit) because of the longer URL text now in the column.
# The COM layer quantizes column widths to whole pixels, so we pick the input
# value that rounds to the closest achievable stored width (51.1666... vs the
# target 51.21484375).
$elements.Columns.Item(26).ColumnWidth = 50.35
